$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1105
$ws1.Range("F4").Value = 247
$ws1.Range("F6").Value = 12050
$ws1.Range("F7").Value = 45
$ws1.Range("F8").Value = 82
$ws1.Range("F9").Value = 11790
$ws1.Range("F10").Value = 4748
$ws1.Range("F11").Value = 540
$ws1.Range("F12").Value = 71
$ws1.Range("F13").Value = 15
$ws1.Range("F16").Value = 930
$ws1.Range("F17").Value = 354
$ws1.Range("F19").Value = 55

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1105
$ws4.Range("F4").Value = 247
$ws4.Range("F8").Value = 12050
$ws4.Range("F9").Value = 45
$ws4.Range("F10").Value = 82
$ws4.Range("F11").Value = 11790
$ws4.Range("F12").Value = 4748
$ws4.Range("F13").Value = 540
$ws4.Range("F14").Value = 71
$ws4.Range("F15").Value = 15
$ws4.Range("F18").Value = 930
$ws4.Range("F19").Value = 354
$ws4.Range("F21").Value = 55
